$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (H1) to the new header cells so the
# new columns I and J get the same bold/centered/bordered formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# New header labels (I0, IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2 and 3
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 5
